$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D (Price) to text format first so purely-numeric-looking
# strings like "602.56" are stored as text, matching the source data,
# instead of being auto-converted to numbers by Excel.
$priceRange = $ws.Range("D2:D51")
$priceRange.NumberFormat = "@"

$ws.Range('D2').Value = '68.935.06'
$ws.Range('E2').Value = '  +3.03%  '
$ws.Range('D3').Value = '3.750.55'
$ws.Range('E3').Value = '  +2.46%  '
$ws.Range('E4').Value = '  -0.11%  '
$ws.Range('D5').Value = '602.56'
$ws.Range('E5').Value = '  +1.91%  '
$ws.Range('D6').Value = '169.15'
$ws.Range('E6').Value = '  +3.08%  '
$ws.Range('D7').Value = '3.753.01'
$ws.Range('E7').Value = '  +2.44%  '
$ws.Range('E8').Value = '  -0.05%  '
$ws.Range('E9').Value = '  +2.88%  '
$ws.Range('D10').Value = '0.166'
$ws.Range('E10').Value = '  +5.64%  '
$ws.Range('D11').Value = '6.36'
$ws.Range('E11').Value = '  +4.37%  '
$ws.Range('D12').Value = '0.464'
$ws.Range('E12').Value = '  +1.29%  '
$ws.Range('D13').Value = '38.30'
$ws.Range('E13').Value = '  +2.79%  '
$ws.Range('E14').Value = '  +4.48%  '
$ws.Range('D15').Value = '4.377.62'
$ws.Range('E15').Value = '  +2.52%  '
$ws.Range('D16').Value = '3.741.46'
$ws.Range('E16').Value = '  +2.10%  '
$ws.Range('D17').Value = '68.910.34'
$ws.Range('E17').Value = '  +2.77%  '
$ws.Range('D18').Value = '7.30'
$ws.Range('E18').Value = '  +2.98%  '
$ws.Range('E19').Value = '  +0.88%  '
$ws.Range('E20').Value = '  +1.73%  '
$ws.Range('B21').Value = 'BitcoinCash'
$ws.Range('C21').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range('D21').Value = '497.83'
$ws.Range('E21').Value = '  +2.33%  '
$ws.Range('B22').Value = 'Uniswap'
$ws.Range('C22').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range('D22').Value = '10.75'
$ws.Range('E22').Value = '  +19.60%  '
$ws.Range('D23').Value = '0.729'
$ws.Range('E23').Value = '  +2.73%  '
$ws.Range('D24').Value = '0.0000153'
$ws.Range('E24').Value = '  +11.33%  '
$ws.Range('D25').Value = '85.50'
$ws.Range('E25').Value = '  +0.59%  '
$ws.Range('D26').Value = '2.33'
$ws.Range('E26').Value = '  +2.57%  '
$ws.Range('D27').Value = '12.36'
$ws.Range('E27').Value = '  +2.84%  '
$ws.Range('D28').Value = '10.32'
$ws.Range('E28').Value = '  +4.84%  '
$ws.Range('E29').Value = '  +0.57%  '
$ws.Range('B30').Value = 'ImmutableX'
$ws.Range('C30').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D30').Value = '2.53'
$ws.Range('E30').Value = '  +8.58%  '
$ws.Range('B31').Value = 'PancakeSwap'
$ws.Range('C31').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D31').Value = '2.99'
$ws.Range('E31').Value = '  +3.22%  '
$ws.Range('D32').Value = '7.92'
$ws.Range('E32').Value = '  +3.61%  '
$ws.Range('D33').Value = '31.86'
$ws.Range('E33').Value = '  +1.84%  '
$ws.Range('D34').Value = '3.894.88'
$ws.Range('E34').Value = '  +2.48%  '
$ws.Range('B35').Value = 'Hedera'
$ws.Range('C35').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D35').Value = '0.109'
$ws.Range('E35').Value = '  +2.92%  '
$ws.Range('B36').Value = 'RenzoRestakedETH'
$ws.Range('C36').Value = 'https://coinranking.com/coin/lKlJ_MC5M+renzorestakedeth-ezeth'
$ws.Range('D36').Value = '3.684.32'
$ws.Range('E36').Value = '  +2.41%  '
$ws.Range('D37').Value = '0.998'
$ws.Range('E37').Value = '  -0.26%  '
$ws.Range('E38').Value = '  +2.96%  '
$ws.Range('D39').Value = '5.87'
$ws.Range('E39').Value = '  +2.91%  '
$ws.Range('E40').Value = '  +1.90%  '
$ws.Range('E41').Value = '  +1.60%  '
$ws.Range('D42').Value = '438.60'
$ws.Range('E42').Value = '  +1.12%  '
$ws.Range('D43').Value = '2.96'
$ws.Range('E43').Value = '  +7.96%  '
$ws.Range('D44').Value = '48.92'
$ws.Range('E44').Value = '  +0.94%  '
$ws.Range('E45').Value = '  +3.48%  '
$ws.Range('D46').Value = '8.48'
$ws.Range('E46').Value = '  +2.64%  '
$ws.Range('E47').Value = '  +0.00%  '
$ws.Range('D48').Value = '40.54'
$ws.Range('E48').Value = '  +2.90%  '
$ws.Range('D49').Value = '141.76'
$ws.Range('E49').Value = '  -0.05%  '
$ws.Range('D50').Value = '2.788.83'
$ws.Range('E50').Value = '  +1.78%  '
$ws.Range('D51').Value = '0.0356'
$ws.Range('E51').Value = '  +3.55%  '

# Restore the default (unstyled) cell style now that the text values are
# committed, so the style index matches the original workbook.
$priceRange.Style = "Normal"
